$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New draw-result row appended for 2025-10-02.
# Force text storage (so date-/number-looking strings like "2025-10-02"
# and "251002" keep their literal text form instead of being coerced
# into a date serial / number), then restore the default "Normal"
# style so no stray per-cell formatting is introduced.
$ws.Range("A16:E16").NumberFormat = "@"

$ws.Range("A16").Value = "2025-10-02"
$ws.Range("B16").Value = "Pick 3"
$ws.Range("C16").Value = "251002"
$ws.Range("D16").Value = "0-4-5"
$ws.Range("E16").Value = "2025-10-02T21:36:12.368+04:00"

$ws.Range("A16:E16").Style = "Normal"
